$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Detect Loop and Remove Loop"
$ws.Range("H15").Value = "RemoveLoop"

$ws.Range("H15").Select()
